$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 2.4
$ws.Range("I4").Value = 3.5
$ws.Range("J4").Value = 3.25
$ws.Range("K4").Value = 1.83
$ws.Range("L4").Value = 4.33
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.78
$ws.Range("S4").Value = 2.88
$ws.Range("T4").Value = 1.4
$ws.Range("U4").Value = 4.7
$ws.Range("V4").Value = 1.19
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 1.13
$ws.Range("AA4").Value = 2.2
$ws.Range("AB4").Value = 1.62
$ws.Range("AI4").Value = 5.5
$ws.Range("AL4").Value = 81
$ws.Range("AQ4").Value = 41
$ws.Range("AS4").Value = 51

# Row 5
$ws.Range("G5").Value = 3.1
$ws.Range("H5").Value = 2.7
$ws.Range("I5").Value = 2.75
$ws.Range("J5").Value = 4.33
$ws.Range("K5").Value = 1.67
$ws.Range("M5").Value = 1.2
$ws.Range("N5").Value = 4.33
$ws.Range("O5").Value = 1.91
$ws.Range("P5").Value = 1.8
$ws.Range("U5").Value = 7.8
$ws.Range("V5").Value = 1.08
$ws.Range("Y5").Value = 1.93
$ws.Range("Z5").Value = 1.88
$ws.Range("AA5").Value = 3
$ws.Range("AB5").Value = 1.36
$ws.Range("AD5").Value = 12
$ws.Range("AE5").Value = 15
$ws.Range("AF5").Value = 41
$ws.Range("AH5").Value = 67
$ws.Range("AI5").Value = 4.33
$ws.Range("AJ5").Value = 6.5
$ws.Range("AN5").Value = 5
$ws.Range("AO5").Value = 11
$ws.Range("AS5").Value = 67

# Row 7
$ws.Range("G7").Value = 2.45
$ws.Range("H7").Value = 2.7
$ws.Range("I7").Value = 3.5
$ws.Range("J7").Value = 3.4
$ws.Range("L7").Value = 4.5
$ws.Range("M7").Value = 1.17
$ws.Range("N7").Value = 5
$ws.Range("AC7").Value = 5.5
$ws.Range("AD7").Value = 9.5
$ws.Range("AE7").Value = 12
$ws.Range("AF7").Value = 23
$ws.Range("AG7").Value = 29
$ws.Range("AI7").Value = 4.75
$ws.Range("AN7").Value = 6.5

# Row 8
$ws.Range("G8").Value = 2.63
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 2.6
$ws.Range("J8").Value = 3.4
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 8.5
$ws.Range("O8").Value = 1.36
$ws.Range("P8").Value = 3
$ws.Range("S8").Value = 2.15
$ws.Range("T8").Value = 1.67
$ws.Range("AD8").Value = 13
$ws.Range("AE8").Value = 11
$ws.Range("AG8").Value = 23
$ws.Range("AI8").Value = 8.5
$ws.Range("AJ8").Value = 6
$ws.Range("AM8").Value = 301
$ws.Range("AN8").Value = 8
$ws.Range("AO8").Value = 12
$ws.Range("AP8").Value = 10
$ws.Range("AQ8").Value = 26

# Row 14
$ws.Range("S14").Value = 2.2
$ws.Range("T14").Value = 1.65
$ws.Range("W14").Value = 4
$ws.Range("X14").Value = 1.22
